# Update Ixion_Profits market-price snapshot values (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 15366.598
$ws.Range("I15").Value = 15366.598
$ws.Range("K15").Value = 46099.794
$ws.Range("M15").Value = -45930.794
$ws.Range("H40").Value = 1960
$ws.Range("I40").Value = 1960
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1960
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1785
$ws.Range("N40").ClearContents()
$ws.Range("H76").Value = 3523.0557
$ws.Range("J76").Value = 4046.889
$ws.Range("L76").Value = 4046.889
$ws.Range("N76").Value = -4676.889
$ws.Range("H79").Value = 3523.0557
$ws.Range("J79").Value = 4046.889
$ws.Range("L79").Value = 4046.889
$ws.Range("N79").Value = -6230.889
$ws.Range("H98").Value = 1122.5
$ws.Range("I98").Value = 1122.5
$ws.Range("K98").Value = 1122.5
$ws.Range("M98").Value = 375.5
$ws.Range("H107").Value = 12503593
$ws.Range("I107").Value = 13159045
$ws.Range("J107").Value = 50000
$ws.Range("K107").Value = 13159045
$ws.Range("L107").Value = 50000
$ws.Range("M107").Value = -13157125
$ws.Range("N107").Value = -53840
$ws.Range("H113").Value = 4891.6772
$ws.Range("I113").Value = 3293.5
$ws.Range("J113").Value = 6207.8237
$ws.Range("K113").Value = 3293.5
$ws.Range("L113").Value = 6207.8237
$ws.Range("M113").Value = -39.5
$ws.Range("N113").Value = -12715.8237
$ws.Range("H122").Value = 1122.5
$ws.Range("I122").Value = 1122.5
$ws.Range("K122").Value = 3367.5
$ws.Range("M122").Value = -917.5
$ws.Range("H137").Value = 1447.4054
$ws.Range("I137").Value = 844.6786
$ws.Range("J137").Value = 3322.5557
$ws.Range("K137").Value = 2534.0358
$ws.Range("L137").Value = 9967.667099999999
$ws.Range("M137").Value = 15.96420000000035
$ws.Range("N137").Value = -15067.6671
$ws.Range("H138").Value = 2387.9895
$ws.Range("I138").Value = 992.1163
$ws.Range("J138").Value = 3520.4905
$ws.Range("K138").Value = 2976.3489
$ws.Range("L138").Value = 10561.4715
$ws.Range("M138").Value = 2163.6511
$ws.Range("N138").Value = -20841.4715

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 125001900
$ws.Range("I63").Value = 142858880
$ws.Range("K63").Value = 142858880
$ws.Range("M63").Value = -142858194
$ws.Range("H66").Value = 125001900
$ws.Range("I66").Value = 142858880
$ws.Range("K66").Value = 714294400
$ws.Range("M66").Value = -714290968
$ws.Range("H88").Value = 125002536
$ws.Range("I88").Value = 2743
$ws.Range("J88").Value = 166669140
$ws.Range("K88").Value = 2743
$ws.Range("L88").Value = 166669140
$ws.Range("M88").Value = -2337
$ws.Range("N88").Value = -166669952
$ws.Range("H91").Value = 125002536
$ws.Range("I91").Value = 2743
$ws.Range("J91").Value = 166669140
$ws.Range("K91").Value = 2743
$ws.Range("L91").Value = 166669140
$ws.Range("M91").Value = -1339
$ws.Range("N91").Value = -166671948
$ws.Range("H132").Value = 3525.8108
$ws.Range("I132").Value = 2678.4285
$ws.Range("J132").Value = 4041.6086
$ws.Range("K132").Value = 8035.2855
$ws.Range("L132").Value = 12124.8258
$ws.Range("M132").Value = -5505.2855
$ws.Range("N132").Value = -17184.8258

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2528.2666
$ws.Range("I86").Value = 2174.889
$ws.Range("J86").Value = 3058.3333
$ws.Range("K86").Value = 2174.889
$ws.Range("L86").Value = 3058.3333
$ws.Range("M86").Value = -1051.889
$ws.Range("N86").Value = -5304.3333
$ws.Range("H89").Value = 2528.2666
$ws.Range("I89").Value = 2174.889
$ws.Range("J89").Value = 3058.3333
$ws.Range("K89").Value = 10874.445
$ws.Range("L89").Value = 15291.6665
$ws.Range("M89").Value = -5258.445
$ws.Range("N89").Value = -26523.6665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2516.9856
$ws.Range("I31").Value = 1649.8823
$ws.Range("J31").Value = 2795.1133
$ws.Range("K31").Value = 1649.8823
$ws.Range("L31").Value = 2795.1133
$ws.Range("M31").Value = -1354.8823
$ws.Range("N31").Value = -3385.1133
$ws.Range("H34").Value = 2516.9856
$ws.Range("I34").Value = 1649.8823
$ws.Range("J34").Value = 2795.1133
$ws.Range("K34").Value = 1649.8823
$ws.Range("L34").Value = 2795.1133
$ws.Range("M34").Value = -1447.8823
$ws.Range("N34").Value = -3199.1133
$ws.Range("H134").Value = 1940.6171
$ws.Range("I134").Value = 2197.5405
$ws.Range("J134").Value = 990
$ws.Range("K134").Value = 6592.6215
$ws.Range("L134").Value = 2970
$ws.Range("M134").Value = -4057.6215
$ws.Range("N134").Value = -8040

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3731
$ws.Range("I68").Value = 4155.8
$ws.Range("J68").Value = 2934.5
$ws.Range("K68").Value = 12467.4
$ws.Range("L68").Value = 8803.5
$ws.Range("M68").Value = -11656.4
$ws.Range("N68").Value = -10425.5
$ws.Range("H71").Value = 3731
$ws.Range("I71").Value = 4155.8
$ws.Range("J71").Value = 2934.5
$ws.Range("K71").Value = 37402.2
$ws.Range("L71").Value = 26410.5
$ws.Range("M71").Value = -33346.2
$ws.Range("N71").Value = -34522.5
$ws.Range("H113").Value = 238595.67
$ws.Range("I113").Value = 476.39285
$ws.Range("J113").Value = 714834.2
$ws.Range("K113").Value = 1429.17855
$ws.Range("L113").Value = 2144502.6
$ws.Range("M113").Value = 740.8214499999999
$ws.Range("N113").Value = -2148842.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5950.1665
$ws.Range("I70").Value = 5940.7715
$ws.Range("K70").Value = 5940.7715
$ws.Range("M70").Value = -5670.7715
$ws.Range("H73").Value = 5950.1665
$ws.Range("I73").Value = 5940.7715
$ws.Range("K73").Value = 5940.7715
$ws.Range("M73").Value = -5004.7715
$ws.Range("H126").Value = 5609.7915
$ws.Range("I126").Value = 6315.952
$ws.Range("J126").Value = 666.6667
$ws.Range("K126").Value = 18947.856
$ws.Range("L126").Value = 2000.0001
$ws.Range("M126").Value = -16477.856
$ws.Range("N126").Value = -6940.0001
$ws.Range("H135").Value = 32976.19
$ws.Range("J135").Value = 32976.19
$ws.Range("L135").Value = 32976.19
$ws.Range("N135").Value = -43116.19

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3176904.8
$ws.Range("I22").Value = 12346813
$ws.Range("J22").Value = 2705.8462
$ws.Range("K22").Value = 12346813
$ws.Range("L22").Value = 2705.8462
$ws.Range("M22").Value = -12346518
$ws.Range("N22").Value = -3295.8462
$ws.Range("H27").Value = 3176904.8
$ws.Range("I27").Value = 12346813
$ws.Range("J27").Value = 2705.8462
$ws.Range("K27").Value = 12346813
$ws.Range("L27").Value = 2705.8462
$ws.Range("M27").Value = -12346706
$ws.Range("N27").Value = -2919.8462
$ws.Range("H136").Value = 3480.9434
$ws.Range("I136").Value = 1784.375
$ws.Range("J136").Value = 8701.154
$ws.Range("K136").Value = 5353.125
$ws.Range("L136").Value = 26103.462
$ws.Range("M136").Value = -2803.125
$ws.Range("N136").Value = -31203.462
